# Update "想去人数" (want-to-go count) values in column F
# for the "展览" and "全部类型" sheets, as produced by the
# gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Row -> new value mapping (F column) shared by both sheets
$updates = @{
    6  = 123
    7  = 1203
    19 = 1701
    23 = 647
    25 = 329
    26 = 4090
    35 = 187
    36 = 47
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
